$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 6 fresh rows right before the current totals row (row 31), which
#    pushes the "tunnit yht." total row down to row 37 and carries its
#    column-B / column-C direct formatting down with it (Excel's normal
#    insert-shift behaviour, formatting inherited from the row above the
#    insertion point).
# ---------------------------------------------------------------------------
$ws.Range("A31:C36").Insert(-4121) | Out-Null

# ---------------------------------------------------------------------------
# 2) New data rows 31-36.
# ---------------------------------------------------------------------------

# Row 31 -- new date-group header "3.21.2021" (kept as literal TEXT, not an
# actual Excel date, same as the source diff: t="s"). Writing the literal
# string directly gets auto-converted to a real date serial by COM's usual
# text-to-date inference, so route it through a TEXT-returning formula and
# then freeze it to a plain value -- this keeps the result a shared-string
# cell with no left-over/alternate number-format style.
$ws.Cells.Item(31,1).Formula = '=T("3.21.2021")'
$ws.Cells.Item(31,1).Copy() | Out-Null
$ws.Cells.Item(31,1).PasteSpecial(-4163) | Out-Null   # xlPasteValues

# Apply the same direct formatting as the old totals-row label cell (center
# horizontal alignment only, General format) by copying formats from A37,
# which is exactly that cell after the shift above.
$ws.Cells.Item(37,1).Copy() | Out-Null
$ws.Cells.Item(31,1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Cells.Item(31,2).Value = 1.5
$ws.Cells.Item(31,3).Value = "refaktorointia, error viestien resetointi Reset napilla, nimien lyhentämistä, DataView komponenttien useEffect"

$ws.Cells.Item(32,2).Value = 1
$ws.Cells.Item(32,3).Value = "DataView css ja näkymän parantelua"

$ws.Cells.Item(33,2).Value = 1
$ws.Cells.Item(33,3).Value = "datan parserit jotta luvut ja päivät näkyy halutulla/siistimmällä tavalla"

$ws.Cells.Item(34,2).Value = 1
$ws.Cells.Item(34,3).Value = "parserit yhteiseen käyttöön utils/functions.js, css väritystä"

$ws.Cells.Item(35,2).Value = 1
$ws.Cells.Item(36,2).Value = 1

# NB: the shared-string table records new unique strings in the order they
# are first written, not in row order. The source workbook's table has
# "kaikki perustoiminnot ..." (row 36) interned one slot before "ajax
# funktio ..." (row 35), so replicate that exact write order here to keep
# sharedStrings.xml index-for-index identical to the target.
$ws.Cells.Item(36,3).Value = "kaikki perustoiminnot saatu toimimaan, date input, fetch, datan manipulointi, datan esittely, perus css"
$ws.Cells.Item(35,3).Value = "ajax funktio luotu, axios implementation myöhemmin"

# ---------------------------------------------------------------------------
# 3) Fix up the totals row (now row 37) so the SUM formula covers the new
#    data range.
# ---------------------------------------------------------------------------
$ws.Cells.Item(37,2).Formula = "=SUM(B2:B36)"

# ---------------------------------------------------------------------------
# 4) View state: active sheet scrolled to A19, selection on C35.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("C35").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5) Workbook window size/position (restored, non-maximized geometry).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.WindowState = -4143   # xlNormal
$excel.ActiveWindow.Left = 3345
$excel.ActiveWindow.Top = 2640
$excel.ActiveWindow.Width = 21600
$excel.ActiveWindow.Height = 11385

$wb.Save()
